$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B33").Value = "GP22"
$ws.Range("B34").Value = "GND"
$ws.Range("B35").Value = "GP21"
$ws.Range("B36").Value = "GP20"
$ws.Range("B37").Value = "GP19"
$ws.Range("B38").Value = "GP18"
$ws.Range("B39").Value = "GND"
$ws.Range("B40").Value = "GP17"
$ws.Range("B41").Value = "GP16"
$ws.Range("A41").Value = ""
$ws.Range("C41").Value = ""

$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C33").Select()
